$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("db")

# Add the new "userid" entry under the social media TODO column (col I) row 6
$ws.Range("I6").Value = "userid"

# Add new TODO items for conditions/search note in col I rows 13-14
$ws.Range("I13").Value = "conditions"
$ws.Range("I14").Value = "lalabas lang sa search pag wala pang naaadd"

# Update the active selection to I14 to match the saved cursor position
$ws.Activate()
$ws.Range("I14").Select()
